# Trade #37 (Trade #9 for the MarketMaking-only sheet's numbering base) is
# closed at 2026-02-17 20:48:47, and a brand-new OPEN trade (#70) is logged.
# This updates the Summary roll-up, the Strategy Status roll-up for
# MarketMaking, the "All Trades" log, and the per-strategy "MarketMaking" log.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.29   # Current Capital
$summary.Range("B4").Value = 0.08      # Total P&L $
$summary.Range("B6").Value = 37        # Total Trades
$summary.Range("B7").Value = 15        # Winning Trades
$summary.Range("B9").Value = 40.54     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.29   # Capital
$status.Range("D5").Value = 4        # Trades
$status.Range("E5").Value = -0.03    # P&L $
$status.Range("F5").Value = 0.29     # P&L %
$status.Range("G5").Value = 25       # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close out existing Trade #37 (row 38)
$allTrades.Cells.Item(38, 7).Value = 0.93           # Exit Price
$allTrades.Cells.Item(38, 8).Value = "CLOSED"       # Status
$allTrades.Cells.Item(38, 9).Value = 1.087          # P&L %
$allTrades.Cells.Item(38, 10).Value = 0.01          # P&L $
$allTrades.Cells.Item(38, 11).Value = 100.29        # Capital After
$allTrades.Cells.Item(38, 12).Value = "early_exit"  # Exit Reason
$allTrades.Cells.Item(38, 13).Value = 0.11          # Duration (min)

# Append new Trade #70 (row 71, a fresh OPEN trade)
$allTrades.Cells.Item(71, 1).Value = 70
$allTrades.Cells.Item(71, 2).NumberFormat = "@"
$allTrades.Cells.Item(71, 2).Value = "'2026-02-17"
$allTrades.Cells.Item(71, 3).Value = "20:48:41"
$allTrades.Cells.Item(71, 4).Value = "MarketMaking"
$allTrades.Cells.Item(71, 5).Value = "DOWN"
$allTrades.Cells.Item(71, 6).Value = 0.92
$allTrades.Cells.Item(71, 8).Value = "OPEN"
$allTrades.Cells.Item(71, 9).Value = 0
$allTrades.Cells.Item(71, 10).Value = 0
$allTrades.Cells.Item(71, 11).Value = 100.2838583996649
$allTrades.Cells.Item(71, 13).Value = 0
$allTrades.Cells.Item(71, 14).Value = 0
$allTrades.Cells.Item(71, 15).Value = 0
$allTrades.Cells.Item(71, 16).Value = 0.6
$allTrades.Cells.Item(71, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet (per-strategy log, different column layout:
# G=Exit Price, H=Status, I=P&L%, J=P&L$, K=Capital After,
# L=Entry Slippage, M=Exit Slippage, N=Confidence, O=Entry Reason,
# P=Exit Reason, Q=Duration (min))
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Close out existing Trade #37 (row 5)
$mm.Cells.Item(5, 7).Value = 0.93            # Exit Price
$mm.Cells.Item(5, 8).Value = "CLOSED"        # Status
$mm.Cells.Item(5, 9).Value = 1.087           # P&L %
$mm.Cells.Item(5, 10).Value = 0.01           # P&L $
$mm.Cells.Item(5, 11).Value = 100.29         # Capital After
$mm.Cells.Item(5, 16).Value = "early_exit"   # Exit Reason
$mm.Cells.Item(5, 17).Value = 0.11           # Duration (min)

# Append new Trade #70 (row 38, a fresh OPEN trade)
$mm.Cells.Item(38, 1).Value = 70
$mm.Cells.Item(38, 2).NumberFormat = "@"
$mm.Cells.Item(38, 2).Value = "'2026-02-17"
$mm.Cells.Item(38, 3).Value = "20:48:41"
$mm.Cells.Item(38, 4).Value = "MarketMaking"
$mm.Cells.Item(38, 5).Value = "DOWN"
$mm.Cells.Item(38, 6).Value = 0.92
$mm.Cells.Item(38, 8).Value = "OPEN"
$mm.Cells.Item(38, 9).Value = 0
$mm.Cells.Item(38, 10).Value = 0
$mm.Cells.Item(38, 11).Value = 100.2838583996649
$mm.Cells.Item(38, 12).Value = 0
$mm.Cells.Item(38, 13).Value = 0
$mm.Cells.Item(38, 14).Value = 0.6
$mm.Cells.Item(38, 15).Value = "Normal spread capture: 19600 bps"
$mm.Cells.Item(38, 17).Value = 0
